# Cargue de plan de desarrollo Ok
#
# The single shared "COD" header (reused for columns A,C,E,G,I) is replaced
# with five distinct headers, one per code column: COD_D, COD_E, COD_P,
# COD_S, COD_M. The other headers (DIMENSION, EJE, PROGRAMA, SUBPROGRAMA,
# META) keep their text/positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "COD_D"
$ws.Range("B1").Value = "DIMENSION"
$ws.Range("C1").Value = "COD_E"
$ws.Range("D1").Value = "EJE"
$ws.Range("E1").Value = "COD_P"
$ws.Range("F1").Value = "PROGRAMA"
$ws.Range("G1").Value = "COD_S"
$ws.Range("H1").Value = "SUBPROGRAMA"
$ws.Range("I1").Value = "COD_M"
$ws.Range("J1").Value = "META"

# Widen the narrow "code" columns slightly (A,C,E,G -> 7.6640625, I -> 7.83203125)
$ws.Range("A:A").ColumnWidth = 7.6640625
$ws.Range("C:C").ColumnWidth = 7.6640625
$ws.Range("E:E").ColumnWidth = 7.6640625
$ws.Range("G:G").ColumnWidth = 7.6640625
$ws.Range("I:I").ColumnWidth = 7.83203125

# Move the active selection
$ws.Range("I13").Select()
